# Apply the edits described by the diff between the original
# tests/test_data/test_global_ar6.xlsx and the updated version:
#  1. Update the "2020" column (G) through "2100" column (O) values for
#     rows 3, 4, 10 and 11 on the "data" sheet.
#  2. Update the sheet view: scroll so column E is the left-most visible
#     column, and change the selection to F2:O16 (active cell F2).
#  3. Nudge the workbook window position (best effort / cosmetic).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$ws.Activate()

# --- Row 3 (AR6 climate diagnostics|Emissions|C2F6|Harmonized) ---
$ws.Range("G3").Value  = 1.6087177454461821
$ws.Range("H3").Value  = 0.84198382161913488
$ws.Range("I3").Value  = 0.70356873633854655
$ws.Range("J3").Value  = 0.54352160754369827
$ws.Range("K3").Value  = 0.47048532309107638
$ws.Range("L3").Value  = 0.4645704816927323
$ws.Range("M3").Value  = 0.4600073568537259
$ws.Range("N3").Value  = 0.45471034728610848
$ws.Range("O3").Value  = 0.44882962125114989

# --- Row 4 (AR6 climate diagnostics|Emissions|CF4|Harmonized) ---
$ws.Range("G4").Value  = 11.1380649
$ws.Range("H4").Value  = 5.8295312999999984
$ws.Range("I4").Value  = 4.8712052000000003
$ws.Range("J4").Value  = 3.7631082
$ws.Range("K4").Value  = 3.2574366000000001
$ws.Range("L4").Value  = 3.2164847999999999
$ws.Range("M4").Value  = 3.1848917000000001
$ws.Range("N4").Value  = 3.1482174999999999
$ws.Range("O4").Value  = 3.1075018999999999

# --- Row 10 (AR6 climate diagnostics|Emissions|NOx|Harmonized) ---
$ws.Range("G10").Value = 1229.7630355010001
$ws.Range("H10").Value = 752.09889483300003
$ws.Range("I10").Value = 911.79343712600019
$ws.Range("J10").Value = 938.14800622400003
$ws.Range("K10").Value = 1014.501330986
$ws.Range("L10").Value = 1040.6988646340001
$ws.Range("M10").Value = 1015.725827259
$ws.Range("N10").Value = 982.03410548100032
$ws.Range("O10").Value = 936.62987589600016

# --- Row 11 (AR6 climate diagnostics|Emissions|Sulfur|Harmonized) ---
$ws.Range("G11").Value = 10715.028417993401
$ws.Range("H11").Value = 10018.589616762691
$ws.Range("I11").Value = 9197.7736844025676
$ws.Range("J11").Value = 8111.9367114782908
$ws.Range("K11").Value = 7196.2146728144562
$ws.Range("L11").Value = 6042.0082227535131
$ws.Range("M11").Value = 5323.6645717810206
$ws.Range("N11").Value = 4897.5210224890106
$ws.Range("O11").Value = 4495.4184032751537

# --- Sheet view: scroll right so column E is left-most, select F2:O16 ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("F2:O16").Select()

# --- Workbook window position (best effort; cosmetic host-window state) ---
$win.Left = 19110
$win.Top = 2835
